# Updated cryptos list on Mon Jul 22 21:41:35 UTC 2024 with GitHub Actions
# Refreshes the Price (D) and Volume(1h) (E) columns for each coin row.
# For D-column values that look like plain decimals (e.g. "10.40"), Excel's
# COM layer would otherwise auto-convert the text to a number (dropping the
# trailing zero / introducing float rounding). Forcing the cell to Text
# format before the assignment keeps the exact original string, and
# resetting the style back to "Normal" afterwards avoids leaving a stray
# explicit number-format style on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.862.91'
$ws.Range('E2').Value = '  -0.38%  '

$ws.Range('D3').Value = '3.461.53'
$ws.Range('E3').Value = '  -1.50%  '

$ws.Range('E4').Value = '  +0.01%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '591.84'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.75%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '180.56'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.84%  '

$ws.Range('E7').Value = '  +2.22%  '

$ws.Range('E8').Value = '  +0.05%  '

$ws.Range('D9').Value = '3.454.77'
$ws.Range('E9').Value = '  -1.49%  '

$ws.Range('E10').Value = '  -0.62%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '6.99'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -2.16%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.428'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -1.99%  '

$ws.Range('D13').Value = '4.054.55'
$ws.Range('E13').Value = '  -1.45%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '32.05'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -1.23%  '

$ws.Range('E15').Value = '  -1.03%  '

$ws.Range('D16').Value = '67.820.06'
$ws.Range('E16').Value = '  -0.31%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.0000177'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -3.15%  '

$ws.Range('D18').Value = '3.457.58'

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.18'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -3.63%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '14.08'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -4.91%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '391.35'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -1.73%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '7.90'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -2.58%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.83'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +2.33%  '

$ws.Range('E24').Value = '  -0.13%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.537'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -1.81%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '71.76'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -2.21%  '

$ws.Range('E27').Value = '  -4.76%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.40'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -3.18%  '

$ws.Range('E29').Value = '  -1.89%  '

$ws.Range('E30').Value = '  +0.28%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.10'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -3.46%  '

$ws.Range('E32').Value = '  -1.39%  '

$ws.Range('E33').Value = '  -5.67%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '23.46'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -3.08%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '7.34'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -1.99%  '

$ws.Range('E36').Value = '  -0.14%  '

$ws.Range('E37').Value = '  -7.61%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '162.65'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -1.05%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.887'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +1.13%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.77'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.35%  '

$ws.Range('E41').Value = '  -5.45%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '4.63'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -3.11%  '

$ws.Range('E43').Value = '  -6.90%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '26.04'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -2.83%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0718'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -3.26%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '26.19'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -6.47%  '

$ws.Range('D47').Value = '2.724.73'
$ws.Range('E47').Value = '  -4.32%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '41.29'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -2.57%  '

$ws.Range('E49').Value = '  -3.07%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '329.15'
$ws.Range('D50').Style = 'Normal'

$ws.Range('E51').Value = '  -4.63%  '
